# "Began to refactor + add external factors to evaluations"
#
# Changes applied (per the OOXML diff):
#  - Лист1!I2: 7 -> 5
#  - Лист1!I3: 7 -> 5
#  - Лист1!I4: 6 -> 7
#    (Лист2/"norm data" I2:I4 are formulas `=Лист1!I#/Лист1!I$13` that
#    recalc automatically to 0.5, 0.5, 0.7 respectively.)
#  - The active sheet switches from "norm data" (tab 2) back to "Лист1"
#    (tab 1), i.e. the workbook's activeTab/tabSelected moves to Лист1.
#  - Лист1's sheet view: zoom 132% -> 166%, scroll reset off "H1", and the
#    selected cell moves from N18 to I15.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Лист1")
$ws2 = $wb.Worksheets.Item("norm data")

# --- data edits on Лист1 (drives the recalculated ratios on "norm data") ---
$ws1.Range("I2").Value = 5
$ws1.Range("I3").Value = 5
$ws1.Range("I4").Value = 7

# --- view/selection state ---
# Make Лист1 the active (selected) sheet/tab again.
$ws1.Activate()

# Лист1's sheetView: drop the old "H1" top-left scroll anchor, bump zoom
# to 166%, and move the selection to I15.
$ws1.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 166
$ws1.Range("I15").Select()
